$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Move the existing K2 ("Test") / K5 (long German note) values over to
# the new O2 / O5 position - they keep their original text, just shift
# two columns to the right to make room for the new "LTSD Parameters"
# block.
# ---------------------------------------------------------------------
$ws.Range("O2").Value = $ws.Range("K2").Value()
$ws.Range("O5").Value = $ws.Range("K5").Value()

# New header over K2 ("LTSD Parameters")
$ws.Range("K2").Value = "LTSD Parameters"

# New "Right"/"Left" sub headers on row 3
$ws.Range("K3").Value = "Right"
$ws.Range("M3").Value = "Left"

# New column headers on row 4
$ws.Range("K4").Value = "Threshols"
$ws.Range("L4").Value = "Win"
$ws.Range("M4").Value = "Threshold"
$ws.Range("N4").Value = "Win"

# ---------------------------------------------------------------------
# New parameter values on row 5 (K5:N5). These must stay plain TEXT
# ("4.3", "100.0", "7.0", "400.0") just like the source workbook, not be
# auto-converted to numbers. Writing them straight into Value always
# gets reinterpreted as a number by the engine, so we stage them as
# text formulas in a scratch row far below the used range, copy, and
# paste-special just the values across - this keeps the literal text
# (and its trailing zero) without touching any cell styles.
# ---------------------------------------------------------------------
$ws.Range("A20").Formula = '="4.3"'
$ws.Range("B20").Formula = '="100.0"'
$ws.Range("C20").Formula = '="7.0"'
$ws.Range("D20").Formula = '="400.0"'
$ws.Range("A20:D20").Copy()
$ws.Range("K5").PasteSpecial(-4163) # xlPasteValues
$ws.Rows(20).Delete()

# Update the view: scroll so column B is the left-most visible column,
# and move the active selection to N5
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("N5").Select()
